# Update latest output (run 187)

$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule" ---
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsSchedule.Range("E2").Value = 1160.94052425
$wsSchedule.Range("F2").Value = 19.19544517609127

# --- Sheet "Detailed" ---
$wsDetailed = $wb.Worksheets.Item("Detailed")

$wsDetailed.Range("B13").Value = 58.43713
$wsDetailed.Range("B14").Value = 58.00918

$wsDetailed.Range("C15").Value = "historical"

$wsDetailed.Range("B16").Value = 43.44
$wsDetailed.Range("C16").Value = "historical"

$wsDetailed.Range("B17").Value = 22.16642

$wsDetailed.Range("B19").Value = 33.60787
$wsDetailed.Range("B20").Value = 0.04804
$wsDetailed.Range("B21").Value = -0.05007
$wsDetailed.Range("B22").Value = 0.51
$wsDetailed.Range("B23").Value = 4.02031
$wsDetailed.Range("B24").Value = 23.30665
$wsDetailed.Range("B25").Value = 23.76872
$wsDetailed.Range("B26").Value = 23.75137
$wsDetailed.Range("B27").Value = 34.37273

$wsDetailed.Range("B33").Value = 0
$wsDetailed.Range("B34").Value = 0.00625
$wsDetailed.Range("B35").Value = -4.32098
$wsDetailed.Range("B36").Value = -3.6481
$wsDetailed.Range("B37").Value = 6.00453
$wsDetailed.Range("B38").Value = 35.93074
$wsDetailed.Range("B39").Value = 46.43491

$wsDetailed.Range("B41").Value = 59.73492
$wsDetailed.Range("B42").Value = 59.37697
